# Fix (bug #2): Error indicador 1.4
# Recomputes indicator C1.4 (alternative_info!L), its normalized counterpart
# (alternatives_norm!D) and the resulting TOPSIS evaluation score (result!B),
# plus updates the run date on the info sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet: info ---
$ws = $wb.Worksheets.Item("info")
# C2 holds the run date as plain text (not an Excel date value), so force
# the cell to Text format before assigning to avoid Excel auto-converting it.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "04/02/23"

# --- Sheet: alternative_info (indicator C1.4, column L) ---
$ws = $wb.Worksheets.Item("alternative_info")
$ws.Range("L5").Value = 0.958603951173407
$ws.Range("L7").Value = 0.8510959963766491
$ws.Range("L8").Value = 1.131928037118369
$ws.Range("L10").Value = 1.182560277804791
$ws.Range("L12").Value = 1.023100414119094
$ws.Range("L13").Value = 1.30738877741298
$ws.Range("L14").Value = 0.9389909155463747
$ws.Range("L15").Value = 1.184877940454469
$ws.Range("L18").Value = 1.282366543874901
$ws.Range("L19").Value = 1.547292573912818
$ws.Range("L20").Value = 1.152900358308349
$ws.Range("L21").Value = 1.378578159221931
$ws.Range("L23").Value = 1.28546803852881
$ws.Range("L24").Value = 1.55

# --- Sheet: alternatives_norm (normalized indicator C1.4, column D) ---
$ws = $wb.Worksheets.Item("alternatives_norm")
$ws.Range("D2").Value = 0.006376773535676842
$ws.Range("D3").Value = 0.006405320792817783
$ws.Range("D4").Value = 0.008495410557767199
$ws.Range("D5").Value = 0.000006652146100452816
$ws.Range("D6").Value = 0.009877642888782158
$ws.Range("D7").Value = 0.000007492425722626506
$ws.Range("D8").Value = 0.000005633550302288344
$ws.Range("D9").Value = 0.01062828957875681
$ws.Range("D10").Value = 0.000005392345451949531
$ws.Range("D11").Value = 0.01337066719238917
$ws.Range("D12").Value = 0.000006232793426408052
$ws.Range("D13").Value = 0.000004877488353766505
$ws.Range("D14").Value = 0.000006791091830708884
$ws.Range("D15").Value = 0.000005381797835843734
$ws.Range("D16").Value = 0.1180883988088304
$ws.Range("D17").Value = 0.6871331928299399
$ws.Range("D18").Value = 0.000004972660559599655
$ws.Range("D19").Value = 0.000004121246132236747
$ws.Range("D20").Value = 0.000005531070824744547
$ws.Range("D21").Value = 0.000004625616250351661
$ws.Range("D22").Value = 0.7164914085030158
$ws.Range("D23").Value = 0.0000049606628438424
$ws.Range("D24").Value = 0.000004114047442372156

# --- Sheet: result (final evaluation score, column B) ---
$ws = $wb.Worksheets.Item("result")
$ws.Range("B2").Value = 0.7681420195985504
$ws.Range("B3").Value = 0.765753543979098
$ws.Range("B4").Value = 0.752827513615253
$ws.Range("B5").Value = 0.7375327174716352
$ws.Range("B6").Value = 0.7307311921971286
$ws.Range("B7").Value = 0.7193490658929405
$ws.Range("B8").Value = 0.7010310536318046
$ws.Range("B9").Value = 0.690179211408826
$ws.Range("B10").Value = 0.6884521093021442
$ws.Range("B11").Value = 0.664568824535941
$ws.Range("B12").Value = 0.6616280448072767
$ws.Range("B13").Value = 0.6527679499729255
$ws.Range("B14").Value = 0.6366139187078655
$ws.Range("B15").Value = 0.6300703855101423
$ws.Range("B16").Value = 0.5253621473226203
$ws.Range("B17").Value = 0.520509647572627
$ws.Range("B18").Value = 0.519129068960003
$ws.Range("B19").Value = 0.5189281321773892
$ws.Range("B20").Value = 0.506393556098291
$ws.Range("B21").Value = 0.5009223403762495
$ws.Range("B22").Value = 0.4565863532244553
$ws.Range("B23").Value = 0.4358329597893421
$ws.Range("B24").Value = 0.2903052681303435

Write-Host "Applied indicator 1.4 fix: updated alternative_info, alternatives_norm, result and run date."
